{"js": "const replacements = [\n  [\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"],\n  [\"49\\u00F72=24, 1\", \"65\\u00F77=9, 2\"],\n  [\"58\\u00F74=14, 2\", \"90\\u00F74=22, 2\"],\n  [\"62\\u00F73=20, 2\", \"55\\u00F77=7, 6\"],\n  [\"77\\u00F73=25, 2\", \"37\\u00F74=9, 1\"],\n  [\"28\\u00F76=4, 4\", \"20\\u00F78=2, 4\"],\n  [\"14\\u00F73=4, 2\", \"93\\u00F76=15, 3\"],\n  [\"63\\u00F76=10, 3\", \"23\\u00F79=2, 5\"],\n  [\"23\\u00F72=11, 1\", \"76\\u00F77=10, 6\"],\n  [\"54\\u00F74=13, 2\", \"51\\u00F74=12, 3\"],\n  [\"64\\u00F79=7, 1\", \"63\\u00F77=9, 0\"],\n  [\"15\\u00F75=3, 0\", \"53\\u00F77=7, 4\"],\n  [\"15\\u00F74=3, 3\", \"34\\u00F76=5, 4\"],\n  [\"28\\u00F74=7, 0\", \"96\\u00F79=10, 6\"],\n  [\"28\\u00F72=14, 0\", \"22\\u00F78=2, 6\"],\n  [\"19\\u00F75=3, 4\", \"95\\u00F78=11, 7\"],\n  [\"30\\u00F72=15, 0\", \"78\\u00F77=11, 1\"],\n  [\"83\\u00F79=9, 2\", \"93\\u00F76=15, 3\"],\n  [\"51\\u00F75=10, 1\", \"52\\u00F75=10, 2\"],\n  [\"76\\u00F72=38, 0\", \"66\\u00F72=33, 0\"],\n  [\"29\\u00F74=7, 1\", \"39\\u00F72=19, 1\"],\n  [\"50\\u00F72=25, 0\", \"71\\u00F72=35, 1\"],\n  [\"69\\u00F73=23, 0\", \"66\\u00F72=33, 0\"],\n  [\"38\\u00F77=5, 3\", \"78\\u00F77=11, 1\"],\n  [\"39\\u00F74=9, 3\", \"99\\u00F74=24, 3\"],\n  [\"40\\u00F72=20, 0\", \"81\\u00F74=20, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"),\n    @(\"49\u00f72=24, 1\", \"65\u00f77=9, 2\"),\n    @(\"58\u00f74=14, 2\", \"90\u00f74=22, 2\"),\n    @(\"62\u00f73=20, 2\", \"55\u00f77=7, 6\"),\n    @(\"77\u00f73=25, 2\", \"37\u00f74=9, 1\"),\n    @(\"28\u00f76=4, 4\", \"20\u00f78=2, 4\"),\n    @(\"14\u00f73=4, 2\", \"93\u00f76=15, 3\"),\n    @(\"63\u00f76=10, 3\", \"23\u00f79=2, 5\"),\n    @(\"23\u00f72=11, 1\", \"76\u00f77=10, 6\"),\n    @(\"54\u00f74=13, 2\", \"51\u00f74=12, 3\"),\n    @(\"64\u00f79=7, 1\", \"63\u00f77=9, 0\"),\n    @(\"15\u00f75=3, 0\", \"53\u00f77=7, 4\"),\n    @(\"15\u00f74=3, 3\", \"34\u00f76=5, 4\"),\n    @(\"28\u00f74=7, 0\", \"96\u00f79=10, 6\"),\n    @(\"28\u00f72=14, 0\", \"22\u00f78=2, 6\"),\n    @(\"19\u00f75=3, 4\", \"95\u00f78=11, 7\"),\n    @(\"30\u00f72=15, 0\", \"78\u00f77=11, 1\"),\n    @(\"83\u00f79=9, 2\", \"93\u00f76=15, 3\"),\n    @(\"51\u00f75=10, 1\", \"52\u00f75=10, 2\"),\n    @(\"76\u00f72=38, 0\", \"66\u00f72=33, 0\"),\n    @(\"29\u00f74=7, 1\", \"39\u00f72=19, 1\"),\n    @(\"50\u00f72=25, 0\", \"71\u00f72=35, 1\"),\n    @(\"69\u00f73=23, 0\", \"66\u00f72=33, 0\"),\n    @(\"38\u00f77=5, 3\", \"78\u00f77=11, 1\"),\n    @(\"39\u00f74=9, 3\", \"99\u00f74=24, 3\"),\n    @(\"40\u00f72=20, 0\", \"81\u00f74=20, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
